$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style (bold, bordered, centered header look) from H1 to I1/J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header labels
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Data rows 2-16: column I is always 1, column J mirrors column H
for ($row = 2; $row -le 16; $row++) {
    $hVal = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value2 = 1
    $ws.Cells.Item($row, 10).Value2 = $hVal
}
